$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = 10
$ws.Range("C5").Value = 15
$ws.Range("C8").Value = 13
$ws.Range("C10").Value = 10
$ws.Range("C13").Value = 13
$ws.Range("C15").Value = 12
$ws.Range("C16").Value = 10
$ws.Range("C17").Value = 13
$ws.Range("C18").Value = 12
$ws.Range("C19").Value = 10
$ws.Range("C21").Value = 15
$ws.Range("C22").Value = 10
$ws.Range("C23").Value = 20
$ws.Range("C27").Value = 18
$ws.Range("C28").Value = 15
$ws.Range("C31").Value = 10
$ws.Range("C33").Value = 16
$ws.Range("C40").Value = 14
$ws.Range("C45").Value = 12
$ws.Range("C47").Value = 14
$ws.Range("C49").Value = 15
$ws.Range("C51").Value = 16
$ws.Range("C52").Value = 14
$ws.Range("C53").Value = 16
$ws.Range("C55").Value = 30
$ws.Range("C56").Value = 16
$ws.Range("C59").Value = 22
$ws.Range("B60").Value = 20
$ws.Range("C60").Value = 23
$ws.Range("C61").Value = 17
$ws.Range("C64").Value = 15
$ws.Range("C68").Value = 13
$ws.Range("C69").Value = 12
$ws.Range("C70").Value = 13
$ws.Range("C72").Value = 10
$ws.Range("C78").Value = 14

$wb.Save()
